$d = $word.ActiveDocument

# --- Update the bookmark-style merge ID in the document's first paragraph ---
# Locate the old ID text (without relying on hardcoded character offsets).
$idRange = $d.Content
[void]$idRange.Find.Execute("**ID__AFFARS_pgi_5337_topic_5__ID**", $true, $false, $false, $false, $false, `
                             $true, 1, $false, "", 0)

$idStart = $idRange.Start
$idEnd = $idRange.End

# The run right after the ID text is a lone trailing space in its own <w:r>.
# Replace the ID text together with that trailing space in one go so the
# extra space-only run disappears and only the updated ID text remains.
$fullRange = $d.Range($idStart, $idEnd + 1)
$fullRange.Text = "**ID__AFFARS_AFMC_PGI_5337_102_90__ID**"

# --- Give that same paragraph the same paragraph border / indent as the ---
# --- rest of the body paragraphs                                       ---
$p = $d.Paragraphs(1)
$p.Format.LeftIndent = 11.25

$borders = $p.Format.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5
